$wb = $excel.ActiveWorkbook

# --- General sheet: update selection, no longer the selected tab ---
$wsGeneral = $wb.Worksheets.Item("General")
$wsGeneral.Range("B6").Select()

# --- BESS sheet: update selection and Pmax (F2) value ---
$wsBESS = $wb.Worksheets.Item("BESS")
$wsBESS.Range("F2").Value = 50
$wsBESS.Range("G8").Select()

# --- Generator sheet: rename profile string, bump Phases (C2), becomes active/selected tab ---
$wsGenerator = $wb.Worksheets.Item("Generator")
$wsGenerator.Range("C2").Value = 3
$wsGenerator.Range("I2").Value = "pv_generation"
$wsGenerator.Activate()
$wsGenerator.Range("J5").Select()
